$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.132.68"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.056.59"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.96"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.669"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.28"
$ws.Range("E7").Value = "  +10.44%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.19"
$ws.Range("E12").Value = "  +8.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.822"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.75"
$ws.Range("E15").Value = "  +10.09%  "
$ws.Range("D16").Value = "2.056.82"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.75"
$ws.Range("E17").Value = "  +32.48%  "
$ws.Range("D18").Value = "37.108.34"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.68"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").Value = "0.0₃0914"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.76"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  +13.48%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.10"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.43"
$ws.Range("E27").Value = "  +5.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.37"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.14"
$ws.Range("E30").Value = "  +10.07%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.86"
$ws.Range("E31").Value = "  +6.44%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +5.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0899"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.27"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  +6.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.35"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.22"
$ws.Range("E40").Value = "  +28.35%  "
$ws.Range("E41").Value = "  +10.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.89"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.15"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.06"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.47"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.00"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").Value = "1.296.45"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.89"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "2.239.17"
$ws.Range("E51").Value = "  -0.48%  "
